$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row 2 labels: B2 and F2 become "total" (C2/D2/E2 already correct).
$ws.Range("B2").Value = "total"
$ws.Range("F2").Value = "total"

# Remove the two "section header only" rows that had no data
# (row 8: "grandes regiões e unidades da federação", row 5: "situação do domicílio").
# Delete the lower-numbered row first is fine too, but deleting the higher one
# first avoids needing to recompute the other row's index.
$ws.Rows(8).Delete()
$ws.Rows(5).Delete()
